$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 181.4
$ws.Range("I33").Value = 190.85715
$ws.Range("J33").Value = 159.33333
$ws.Range("K33").Value = 190.85715
$ws.Range("L33").Value = 159.33333
$ws.Range("M33").Value = 38.14285000000001
$ws.Range("N33").Value = -617.3333299999999
$ws.Range("H43").Value = 1829.3334
$ws.Range("I43").Value = 1725
$ws.Range("J43").Value = 1881.5
$ws.Range("K43").Value = 1725
$ws.Range("L43").Value = 1881.5
$ws.Range("M43").Value = -1656
$ws.Range("N43").Value = -2019.5
$ws.Range("H51").Value = 9029.385
$ws.Range("J51").Value = 8784
$ws.Range("L51").Value = 8784
$ws.Range("N51").Value = -9752
$ws.Range("H58").Value = 142.3
$ws.Range("J58").Value = 10
$ws.Range("L58").Value = 30
$ws.Range("N58").Value = -330
$ws.Range("H62").Value = 36209.47
$ws.Range("I62").Value = 6438.857
$ws.Range("J62").Value = 57048.9
$ws.Range("K62").Value = 6438.857
$ws.Range("L62").Value = 57048.9
$ws.Range("M62").Value = -5814.857
$ws.Range("N62").Value = -58296.9
$ws.Range("H65").Value = 36209.47
$ws.Range("I65").Value = 6438.857
$ws.Range("J65").Value = 57048.9
$ws.Range("K65").Value = 32194.285
$ws.Range("L65").Value = 285244.5
$ws.Range("M65").Value = -29074.285
$ws.Range("N65").Value = -291484.5
$ws.Range("H88").Value = 8773207
$ws.Range("I88").Value = 639.2857
$ws.Range("J88").Value = 13890538
$ws.Range("K88").Value = 639.2857
$ws.Range("L88").Value = 13890538
$ws.Range("M88").Value = -233.2857
$ws.Range("N88").Value = -13891350
$ws.Range("H91").Value = 8773207
$ws.Range("I91").Value = 639.2857
$ws.Range("J91").Value = 13890538
$ws.Range("K91").Value = 639.2857
$ws.Range("L91").Value = 13890538
$ws.Range("M91").Value = 764.7143
$ws.Range("N91").Value = -13893346
$ws.Range("H98").Value = 9520.5
$ws.Range("I98").Value = 9453.571
$ws.Range("K98").Value = 9453.571
$ws.Range("M98").Value = -7955.571
$ws.Range("H113").Value = 100003944
$ws.Range("I113").Value = 142860930
$ws.Range("J113").Value = 4332.6665
$ws.Range("K113").Value = 142860930
$ws.Range("L113").Value = 4332.6665
$ws.Range("M113").Value = -142857676
$ws.Range("N113").Value = -10840.6665
$ws.Range("H122").Value = 9520.5
$ws.Range("I122").Value = 9453.571
$ws.Range("K122").Value = 28360.713
$ws.Range("M122").Value = -25910.713
$ws.Range("H132").Value = 10206179
$ws.Range("I132").Value = 11113042
$ws.Range("K132").Value = 33339126
$ws.Range("M132").Value = -33336596
$ws.Range("H137").Value = 4784.48
$ws.Range("I137").Value = 3943.7856
$ws.Range("J137").Value = 5854.4546
$ws.Range("K137").Value = 11831.3568
$ws.Range("L137").Value = 17563.3638
$ws.Range("M137").Value = -9281.356800000001
$ws.Range("N137").Value = -22663.3638
$ws.Range("H138").Value = 2673.1013
$ws.Range("I138").Value = 2324.625
$ws.Range("K138").Value = 6973.875
$ws.Range("M138").Value = -1833.875

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 8248.18
$ws.Range("I32").Value = 5638.4155
$ws.Range("J32").Value = 29363.545
$ws.Range("K32").Value = 5638.4155
$ws.Range("L32").Value = 29363.545
$ws.Range("M32").Value = -5351.4155
$ws.Range("N32").Value = -29937.545
$ws.Range("H81").Value = 52999.5
$ws.Range("I81").Value = 52999.5
$ws.Range("K81").Value = 52999.5
$ws.Range("M81").Value = -52001.5
$ws.Range("H84").Value = 52999.5
$ws.Range("I84").Value = 52999.5
$ws.Range("K84").Value = 158998.5
$ws.Range("M84").Value = -154006.5

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 2271.8572
$ws.Range("I99").Value = 1669.7
$ws.Range("J99").Value = 3777.25
$ws.Range("K99").Value = 1669.7
$ws.Range("L99").Value = 3777.25
$ws.Range("M99").Value = -171.7
$ws.Range("N99").Value = -6773.25

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -324
$ws.Range("H31").Value = 4577.5454
$ws.Range("I31").Value = 4799.6665
$ws.Range("J31").Value = 4494.25
$ws.Range("K31").Value = 4799.6665
$ws.Range("L31").Value = 4494.25
$ws.Range("M31").Value = -4504.6665
$ws.Range("N31").Value = -5084.25
$ws.Range("H34").Value = 4577.5454
$ws.Range("I34").Value = 4799.6665
$ws.Range("J34").Value = 4494.25
$ws.Range("K34").Value = 4799.6665
$ws.Range("L34").Value = 4494.25
$ws.Range("M34").Value = -4597.6665
$ws.Range("N34").Value = -4898.25
$ws.Range("H62").Value = 103600.266
$ws.Range("I62").Value = 4280.8
$ws.Range("K62").Value = 4280.8
$ws.Range("M62").Value = -3656.8
$ws.Range("H65").Value = 103600.266
$ws.Range("I65").Value = 4280.8
$ws.Range("K65").Value = 21404
$ws.Range("M65").Value = -18284
$ws.Range("H108").Value = 78599.2
$ws.Range("J108").Value = 71999
$ws.Range("L108").Value = 71999
$ws.Range("N108").Value = -79679
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178
$ws.Range("H120").Value = 24999.334
$ws.Range("J120").Value = 24999.334
$ws.Range("L120").Value = 24999.334
$ws.Range("N120").Value = -32257.334
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080
$ws.Range("H132").Value = 347720.94
$ws.Range("I132").Value = 2955.5
$ws.Range("K132").Value = 8866.5
$ws.Range("M132").Value = -6336.5
$ws.Range("H133").Value = 74975
$ws.Range("J133").Value = 74975
$ws.Range("L133").Value = 74975
$ws.Range("N133").Value = -80035
$ws.Range("H139").Value = 84825
$ws.Range("J139").Value = 84825
$ws.Range("L139").Value = 84825
$ws.Range("N139").Value = -95105
$ws.Range("H141").Value = 205409.8
$ws.Range("J141").Value = 205409.8
$ws.Range("L141").Value = 205409.8
$ws.Range("N141").Value = -215769.8

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H68").Value = 10667176
$ws.Range("I68").Value = 24381288
$ws.Range("J68").Value = 643.6667
$ws.Range("K68").Value = 73143864
$ws.Range("L68").Value = 1931.0001
$ws.Range("M68").Value = -73143053
$ws.Range("N68").Value = -3553.0001
$ws.Range("H71").Value = 10667176
$ws.Range("I71").Value = 24381288
$ws.Range("J71").Value = 643.6667
$ws.Range("K71").Value = 219431592
$ws.Range("L71").Value = 5793.0003
$ws.Range("M71").Value = -219427536
$ws.Range("N71").Value = -13905.0003
$ws.Range("H131").Value = 9744.414000000001
$ws.Range("J131").Value = 13209.7
$ws.Range("L131").Value = 39629.10000000001
$ws.Range("N131").Value = -49709.10000000001
$ws.Range("H132").Value = 2595.6125
$ws.Range("I132").Value = 1600.5
$ws.Range("J132").Value = 2884.516
$ws.Range("K132").Value = 14404.5
$ws.Range("L132").Value = 25960.644
$ws.Range("M132").Value = -11874.5
$ws.Range("N132").Value = -31020.644

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 11450
$ws.Range("I70").Value = 10501
$ws.Range("J70").Value = 12082.667
$ws.Range("K70").Value = 10501
$ws.Range("L70").Value = 12082.667
$ws.Range("M70").Value = -10231
$ws.Range("N70").Value = -12622.667
$ws.Range("H73").Value = 11450
$ws.Range("I73").Value = 10501
$ws.Range("J73").Value = 12082.667
$ws.Range("K73").Value = 10501
$ws.Range("L73").Value = 12082.667
$ws.Range("M73").Value = -9565
$ws.Range("N73").Value = -13954.667
$ws.Range("H122").Value = 399237.94
$ws.Range("I122").Value = 465152.6
$ws.Range("K122").Value = 1395457.8
$ws.Range("M122").Value = -1393007.8
$ws.Range("H126").Value = 8897.429
$ws.Range("I126").Value = 12986
$ws.Range("K126").Value = 38958
$ws.Range("M126").Value = -36488

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 23813156
$ws.Range("I7").Value = 35716876
$ws.Range("J7").Value = 5713.857
$ws.Range("K7").Value = 35716876
$ws.Range("L7").Value = 5713.857
$ws.Range("M7").Value = -35716764
$ws.Range("N7").Value = -5937.857
$ws.Range("H22").Value = 692.8570999999999
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 692.8570999999999
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 3369.6191
$ws.Range("I40").Value = 2651.2
$ws.Range("J40").Value = 5165.6665
$ws.Range("K40").Value = 2651.2
$ws.Range("L40").Value = 5165.6665
$ws.Range("M40").Value = -2515.2
$ws.Range("N40").Value = -5437.6665
$ws.Range("H122").Value = 8257.454
$ws.Range("I122").Value = 7509.706
$ws.Range("J122").Value = 10799.8
$ws.Range("K122").Value = 22529.118
$ws.Range("L122").Value = 32399.4
$ws.Range("M122").Value = -20079.118
$ws.Range("N122").Value = -37299.39999999999
$ws.Range("H126").Value = 23813156
$ws.Range("I126").Value = 35716876
$ws.Range("J126").Value = 5713.857
$ws.Range("K126").Value = 107150628
$ws.Range("L126").Value = 17141.571
$ws.Range("M126").Value = -107148158
$ws.Range("N126").Value = -22081.571

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H16").Value = 69950
$ws.Range("J16").Value = 69950
$ws.Range("L16").Value = 69950
$ws.Range("N16").Value = -70534
$ws.Range("H62").Value = 1593074.6
$ws.Range("I62").Value = 3404360.2
$ws.Range("K62").Value = 3404360.2
$ws.Range("M62").Value = -3403736.2
$ws.Range("H65").Value = 1593074.6
$ws.Range("I65").Value = 3404360.2
$ws.Range("K65").Value = 17021801
$ws.Range("M65").Value = -17018681
$ws.Range("H132").Value = 487908.1
$ws.Range("I132").Value = 592669
$ws.Range("K132").Value = 1778007
$ws.Range("M132").Value = -1775477
